$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order of writes matters for shared-string table ordering:
# text (C7) first, then image_url (E7), then the repeated "Ehrlichkeit" (B7 & D7)
$ws.Range("C7").Value = "Ich war gestern leicht angeschlagen und habe mir deshalb rausgenommen, erst spät so gegen 10:30 Uhr ins Büro zu gehen. Ich habe diese Zeit einfach gebraucht, um im Tag anzukommen und um es überhaupt ins Büro zu schaffen. Das Gefühl von `"ich muss doch`" ist nicht aufgekommen, sondern ich habe es einfach gemacht. `nWeiter kam die Frage auf, wie ich einen Call mit einem möglichen Geschäftspartner gefunden habe. Ich habe mich hier ehrlich geäussert und gesagt, dass ich mich eigentlich nicht als Programmierer sehe sprich diese Aufgabe für mich mittelfristig nichts ist. Da ich weiss, dass ich dafür vorgesehen bin, hat es schon etwas Überwindung gebraucht. Ich habe bei anderen Geschäften auch schon klar meine Meinung gesagt, dass wir diese Verträge aus meiner Sicht kündigen sollten, weil wir dieses Know-How gar nicht im Unternehmen haben. Er wollte das nicht und jetzt muss ein Arbeitskollege diese Aufträge wahrnehmen, obwohl ihm die Skills dazu eigentlich auch fehlen. Es geht einfach nur darum, Aufträge die Geld geben, nicht zu verlieren statt sich auf das zu fokussieren, was wir können und wollen. Ich bleibe hier mutig und äussere meine Meinung, auch wenn ich damit nicht immer auf Freude stosse. Ich bleibe mir treu."

# The embedded line break makes the engine auto-expand the row height (ht=17.75,
# customHeight=1). Re-run AutoFit so the row reverts to the sheet's standard
# (non-custom) height, matching the target workbook.
$ws.Range("C7").EntireRow.AutoFit()

$ws.Range("E7").Value = "https://previews.123rf.com/images/rashadashurov/rashadashurov1910/rashadashurov191003965/131279365-honesty-icon-simple-element-illustration-honesty-concept-symbol-design-can-be-used-for-web.avif"
$ws.Range("B7").Value = "Ehrlichkeit"
$ws.Range("D7").Value = "Ehrlichkeit"

$ws.Range("B7").Select()

$excel.ActiveWindow.Zoom = 188
